$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (columns A-F)
$data = @(
    @(1, "Literature Society IITJ Website", 87.69, 0.9, 78.92, 3),
    @(1, "Multi Model Data Analysis for Annotation of Human Activities", 78.28, 1, 78.28, 3),
    @(1, "CloudPhysician's Vital Extraction Challenge", 75, 0.85, 63.75, 3),
    @(2, "SMART SENSING MIDDLEWARE", 87.69, 1, 87.69, 3),
    @(2, "RAPID", 81.62, 1, 81.62, 3),
    @(2, "SHAMIYANA APP", 78.28, 0.9, 70.45, 3),
    @(3, "Website for the Literature Society of the college", 87.69, 0.9, 78.92, 2),
    @(3, "LLMGuard", 78.28, 1, 78.28, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Delete rows 10-13 which are no longer part of the data
$ws.Range("A10:F13").Delete() | Out-Null
